$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: header "Save" using the same formatting as the other
# header cells (copy the format from G1, the last existing header).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows: "Save" column values (all 0 for this export)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
